# data update 11 june 20
# Adds a new day (2020-06-11, serial 43993) row to each of the three
# sheets: Confirmed, Recoverd (Recovered), Death.

$wb = $excel.ActiveWorkbook

$confirmed = $wb.Worksheets.Item("Confirmed")
$recovered = $wb.Worksheets.Item("Recoverd")
$death     = $wb.Worksheets.Item("Death")

# --- Confirmed (sheet1): row 97 ---
$confirmed.Range("A96:C96").Copy()
$confirmed.Range("A97:C97").PasteSpecial(-4122)
$confirmed.Cells.Item(97, 1).Value = 43993
$confirmed.Cells.Item(97, 2).Formula = "=SUM(B96+C97)"
$confirmed.Cells.Item(97, 3).Value = 3187

# --- Recoverd (sheet2): row 97 ---
$recovered.Range("A96:C96").Copy()
$recovered.Range("A97:C97").PasteSpecial(-4122)
$recovered.Cells.Item(97, 1).Value = 43993
$recovered.Cells.Item(97, 2).Formula = "=SUM(B96+C97)"
$recovered.Cells.Item(97, 3).Value = 848

# --- Death (sheet3): row 97 ---
$death.Range("A96:C96").Copy()
$death.Range("A97:C97").PasteSpecial(-4122)
$death.Cells.Item(97, 1).Value = 43993
$death.Cells.Item(97, 2).Formula = "=SUM(B96+C97)"
$death.Cells.Item(97, 3).Value = 37

# Mirror the authored selection move (row 96/97 now the "active" pair)
$confirmed.Range("B96:B97").Select()
$recovered.Range("B96:B97").Select()
$death.Range("B96:B97").Select()
